{"js": "// Commit: \"change line spacing and add table numbers\"\n//\n// 1. The paragraph that holds the `{{table:ecosystems}}` merge field used to\n//    carry its own direct formatting (Normal1 + a paragraph border/shading/\n//    line-spacing override + an explicit Calibri/#333333 run) instead of\n//    using the dedicated \"TableCaption\" paragraph style. Re-point it at the\n//    \"TableCaption\" style so the direct formatting is dropped.\n// 2. The \"TableCaption\" style itself switches its line spacing from\n//    312 (15.6pt, auto rule) down to 240 (12pt, auto rule) -- i.e. single\n//    spacing.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find((p) => p.text.indexOf(\"{{table:ecosystems}}\") !== -1);\nif (!target) {\n  throw new Error(\"Could not find the '{{table:ecosystems}}' paragraph\");\n}\n\n// Re-home the paragraph on the TableCaption style; this drops the paragraph's\n// direct pBdr/shd/spacing + the run's direct rFonts/color overrides, since\n// TableCaption already supplies them.\ntarget.style = \"TableCaption\";\n\n// Tighten the TableCaption style's line spacing (15.6pt -> 12pt / single).\nconst styles = context.document.getStyles();\nconst tableCaption = styles.getByNameOrNullObject(\"TableCaption\");\ntableCaption.load(\"nameLocal\");\nawait context.sync();\n\nif (tableCaption.isNullObject) {\n  throw new Error(\"TableCaption style not found\");\n}\ntableCaption.paragraphFormat.lineSpacing = 12;\n\nawait context.sync();\n", "ps1": "# Commit: \"change line spacing and add table numbers\"\n#\n# 1. The paragraph holding the `{{table:ecosystems}}` merge field used to\n#    carry its own direct formatting (Normal1 style + a paragraph border /\n#    shading / line-spacing override + an explicit Calibri/#333333 run)\n#    instead of using the dedicated \"TableCaption\" paragraph style. Re-point\n#    it at the \"TableCaption\" style so that direct formatting is dropped.\n# 2. The \"TableCaption\" style itself switches its line spacing from\n#    312 (15.6pt, auto rule) down to 240 (12pt, auto rule) -- i.e. single\n#    spacing.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*{{table:ecosystems}}*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the '{{table:ecosystems}}' paragraph\"\n}\n\n# Re-home the paragraph on the TableCaption style; this drops the paragraph's\n# direct pBdr/shd/spacing + the run's direct rFonts/color overrides, since\n# TableCaption already supplies them.\n$target.Style = \"TableCaption\"\n\n# Tighten the TableCaption style's line spacing (15.6pt -> 12pt / single).\n$tableCaption = $d.Styles(\"TableCaption\")\n$tableCaption.ParagraphFormat.LineSpacing = 12\n"}
